$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# Per-language, per-row metadata needed to populate the new "Latest Target
# File" (F) / "Latest Handback File" (G) columns, their hyperlinks, and the
# "Latest Handback DateTime" (H) column, for both the zh-cn and de-de sheets.
$sheetsInfo = @(
  @{
    Name = "zh-cn"
    HandbackDateTime = "2016-03-22 20:16:57"
    Rows = @(
      @{
        Row = 2
        MdName = "2fc04ec1-4d65-40a6-b5c6-9171a1d2a072.md"
        MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/38b6a187239a908f6e2c8d8c6c7a30d9c92cd7ca/e2e/2fc04ec1-4d65-40a6-b5c6-9171a1d2a072.md"
        XlfName = "2fc04ec1-4d65-40a6-b5c6-9171a1d2a072.7795558efdc791ff7eea5d6c94f0ba8cd7b0bfc4.zh-cn.xlf"
        XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ecd1841ad887d6b5751cf7bbd18ff78141d86c05/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/2fc04ec1-4d65-40a6-b5c6-9171a1d2a072.7795558efdc791ff7eea5d6c94f0ba8cd7b0bfc4.zh-cn.xlf"
      },
      @{
        Row = 3
        MdName = "a76f29f6-9211-483a-b39b-1a7be9188958.md"
        MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/38b6a187239a908f6e2c8d8c6c7a30d9c92cd7ca/e2e/a76f29f6-9211-483a-b39b-1a7be9188958.md"
        XlfName = "a76f29f6-9211-483a-b39b-1a7be9188958.8d4b360af3b24ae277ff9fdcd5e0ca39fd7f4073.zh-cn.xlf"
        XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ecd1841ad887d6b5751cf7bbd18ff78141d86c05/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/a76f29f6-9211-483a-b39b-1a7be9188958.8d4b360af3b24ae277ff9fdcd5e0ca39fd7f4073.zh-cn.xlf"
      }
    )
  },
  @{
    Name = "de-de"
    HandbackDateTime = "2016-03-22 20:17:03"
    Rows = @(
      @{
        Row = 2
        MdName = "2fc04ec1-4d65-40a6-b5c6-9171a1d2a072.md"
        MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/38b6a187239a908f6e2c8d8c6c7a30d9c92cd7ca/e2e/2fc04ec1-4d65-40a6-b5c6-9171a1d2a072.md"
        XlfName = "2fc04ec1-4d65-40a6-b5c6-9171a1d2a072.7795558efdc791ff7eea5d6c94f0ba8cd7b0bfc4.de-de.xlf"
        XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b80c9fa3b745a9b23ea590170949c964eea96ded/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/2fc04ec1-4d65-40a6-b5c6-9171a1d2a072.7795558efdc791ff7eea5d6c94f0ba8cd7b0bfc4.de-de.xlf"
      },
      @{
        Row = 3
        MdName = "a76f29f6-9211-483a-b39b-1a7be9188958.md"
        MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/38b6a187239a908f6e2c8d8c6c7a30d9c92cd7ca/e2e/a76f29f6-9211-483a-b39b-1a7be9188958.md"
        XlfName = "a76f29f6-9211-483a-b39b-1a7be9188958.8d4b360af3b24ae277ff9fdcd5e0ca39fd7f4073.de-de.xlf"
        XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b80c9fa3b745a9b23ea590170949c964eea96ded/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/a76f29f6-9211-483a-b39b-1a7be9188958.8d4b360af3b24ae277ff9fdcd5e0ca39fd7f4073.de-de.xlf"
      }
    )
  }
)

foreach ($sheetInfo in $sheetsInfo) {
  $ws = $wb.Worksheets.Item($sheetInfo.Name)

  foreach ($rowInfo in $sheetInfo.Rows) {
    $r = $rowInfo.Row

    # Status column now reflects a completed handback.
    $ws.Cells.Item($r, 3).Value = $status

    # Column F: "Latest Target File" - the source .md file, now also the
    # handback target.
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $rowInfo.MdUrl, "", "", $rowInfo.MdName)

    # Column G: "Latest Handback File" - the localized .xlf file that was
    # handed back.
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 7), $rowInfo.XlfUrl, "", "", $rowInfo.XlfName)

    # Column H: "Latest Handback DateTime" - when the handback completed.
    $ws.Cells.Item($r, 8).Value = $sheetInfo.HandbackDateTime
  }
}

# The "Overview" sheet summarises each language's status in its own column
# (B = zh-cn, C = de-de) using the same status text as the per-language
# sheets, so it needs to be kept in sync as well.
$overview = $wb.Worksheets.Item("Overview")
$overview.Cells.Item(2, 2).Value = $status
$overview.Cells.Item(2, 3).Value = $status
$overview.Cells.Item(3, 2).Value = $status
$overview.Cells.Item(3, 3).Value = $status
